$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.020.06'
$ws.Range("E2").Value = '  -5.11%  '
$ws.Range("D3").Value = '3.315.90'
$ws.Range("E3").Value = '  -5.17%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '568.62'
$ws.Range("E5").Value = '  -3.51%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '126.60'
$ws.Range("E6").Value = '  -5.56%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.314.68'
$ws.Range("E8").Value = '  -5.16%  '
$ws.Range("E9").Value = '  -2.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.27'
$ws.Range("E10").Value = '  -4.51%  '
$ws.Range("E11").Value = '  -5.99%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.372'
$ws.Range("E12").Value = '  -4.17%  '
$ws.Range("D13").Value = '3.883.14'
$ws.Range("E13").Value = '  -5.04%  '
$ws.Range("E14").Value = '  -1.47%  '
$ws.Range("D15").Value = '3.322.15'
$ws.Range("E15").Value = '  -5.21%  '
$ws.Range("E16").Value = '  -7.62%  '
$ws.Range("D17").Value = '61.137.16'
$ws.Range("E17").Value = '  -4.88%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '24.34'
$ws.Range("E18").Value = '  -3.86%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.56'
$ws.Range("E19").Value = '  -3.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.05'
$ws.Range("E20").Value = '  -9.82%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.14'
$ws.Range("E21").Value = '  -3.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '349.40'
$ws.Range("E22").Value = '  -9.96%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.551'
$ws.Range("E23").Value = '  -5.28%  '
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("D25").Value = '3.450.19'
$ws.Range("E25").Value = '  -5.09%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '69.98'
$ws.Range("E26").Value = '  -5.59%  '
$ws.Range("E27").Value = '  -8.19%  '
$ws.Range("E28").Value = '  +0.17%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.12'
$ws.Range("E29").Value = '  -3.49%  '
$ws.Range("E30").Value = '  -4.41%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.81'
$ws.Range("E31").Value = '  -4.50%  '
$ws.Range("E32").Value = '  -6.84%  '
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("E34").Value = '  -4.83%  '
$ws.Range("D35").Value = '3.346.57'
$ws.Range("E35").Value = '  -5.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '22.34'
$ws.Range("E36").Value = '  -4.20%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.37'
$ws.Range("E37").Value = '  +0.82%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.66'
$ws.Range("E38").Value = '  -3.87%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '162.67'
$ws.Range("E39").Value = '  -1.73%  '
$ws.Range("E40").Value = '  -4.82%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0749'
$ws.Range("E41").Value = '  -4.76%  '
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.90'
$ws.Range("E43").Value = '  -2.35%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.743'
$ws.Range("E44").Value = '  -7.96%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.20'
$ws.Range("E45").Value = '  -5.05%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.10'
$ws.Range("E46").Value = '  -6.57%  '
$ws.Range("E47").Value = '  -6.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.38'
$ws.Range("E48").Value = '  -8.58%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.61'
$ws.Range("E49").Value = '  -3.08%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.845'
$ws.Range("E50").Value = '  -8.15%  '
$ws.Range("D51").Value = '2.197.19'
$ws.Range("E51").Value = '  -9.18%  '
